$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "298.00"
$ws.Cells.Item(2,4).ClearFormats()
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = "-2.57%"
$ws.Cells.Item(2,5).ClearFormats()

$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "31.65"
$ws.Cells.Item(3,4).ClearFormats()
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = "-2.96%"
$ws.Cells.Item(3,5).ClearFormats()

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "5.170"
$ws.Cells.Item(4,4).ClearFormats()
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = "-2.46%"
$ws.Cells.Item(4,5).ClearFormats()

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "0.07485"
$ws.Cells.Item(5,4).ClearFormats()
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = "0.81%"
$ws.Cells.Item(5,5).ClearFormats()

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "7.782"
$ws.Cells.Item(6,4).ClearFormats()
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = "0.33%"
$ws.Cells.Item(6,5).ClearFormats()

$ws.Cells.Item(7,2).Value = "GateToken"
$ws.Cells.Item(7,3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "3.797"
$ws.Cells.Item(7,4).ClearFormats()
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = "2.48%"
$ws.Cells.Item(7,5).ClearFormats()

$ws.Cells.Item(8,2).Value = "FTXToken"
$ws.Cells.Item(8,3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "1.665"
$ws.Cells.Item(8,4).ClearFormats()
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = "5.87%"
$ws.Cells.Item(8,5).ClearFormats()

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.9246"
$ws.Cells.Item(9,4).ClearFormats()
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = "0.28%"
$ws.Cells.Item(9,5).ClearFormats()

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.1707"
$ws.Cells.Item(10,4).ClearFormats()
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = "2.37%"
$ws.Cells.Item(10,5).ClearFormats()

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.07514"
$ws.Cells.Item(11,4).ClearFormats()
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = "1.16%"
$ws.Cells.Item(11,5).ClearFormats()

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "0.07992"
$ws.Cells.Item(12,4).ClearFormats()
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = "-0.15%"
$ws.Cells.Item(12,5).ClearFormats()

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "0.02988"
$ws.Cells.Item(13,4).ClearFormats()
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = "-3.81%"
$ws.Cells.Item(13,5).ClearFormats()

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "0.09890"
$ws.Cells.Item(14,4).ClearFormats()
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = "0.33%"
$ws.Cells.Item(14,5).ClearFormats()

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "0.001498"
$ws.Cells.Item(15,4).ClearFormats()
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = "-1.89%"
$ws.Cells.Item(15,5).ClearFormats()

$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "0.04665"
$ws.Cells.Item(16,4).ClearFormats()

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "0.006548"
$ws.Cells.Item(17,4).ClearFormats()
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = "7.65%"
$ws.Cells.Item(17,5).ClearFormats()

$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = "-0.88%"
$ws.Cells.Item(18,5).ClearFormats()

$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = "-0.58%"
$ws.Cells.Item(19,5).ClearFormats()

$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = "0.62%"
$ws.Cells.Item(20,5).ClearFormats()

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "0.1337"
$ws.Cells.Item(21,4).ClearFormats()
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value = "1.31%"
$ws.Cells.Item(21,5).ClearFormats()

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "4.568"
$ws.Cells.Item(22,4).ClearFormats()
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = "7.53%"
$ws.Cells.Item(22,5).ClearFormats()

$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = "-4.21%"
$ws.Cells.Item(23,5).ClearFormats()

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "0.001220"
$ws.Cells.Item(24,4).ClearFormats()
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value = "-0.41%"
$ws.Cells.Item(24,5).ClearFormats()

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "0.004413"
$ws.Cells.Item(25,4).ClearFormats()
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value = "-2.72%"
$ws.Cells.Item(25,5).ClearFormats()

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "0.0001400"
$ws.Cells.Item(26,4).ClearFormats()
$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,5).Value = "19.75%"
$ws.Cells.Item(26,5).ClearFormats()

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "0.0001811"
$ws.Cells.Item(27,4).ClearFormats()
$ws.Cells.Item(27,5).NumberFormat = "@"
$ws.Cells.Item(27,5).Value = "8.80%"
$ws.Cells.Item(27,5).ClearFormats()

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "0.01657"
$ws.Cells.Item(39,4).ClearFormats()
$ws.Cells.Item(39,5).NumberFormat = "@"
$ws.Cells.Item(39,5).Value = "2.68%"
$ws.Cells.Item(39,5).ClearFormats()

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "0.04530"
$ws.Cells.Item(40,4).ClearFormats()
$ws.Cells.Item(40,5).NumberFormat = "@"
$ws.Cells.Item(40,5).Value = "0.84%"
$ws.Cells.Item(40,5).ClearFormats()

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "0.006985"
$ws.Cells.Item(41,4).ClearFormats()
$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,5).Value = "-3.89%"
$ws.Cells.Item(41,5).ClearFormats()

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.1342"
$ws.Cells.Item(42,4).ClearFormats()
$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value = "-1.77%"
$ws.Cells.Item(42,5).ClearFormats()

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "0.002060"
$ws.Cells.Item(43,4).ClearFormats()
$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,5).Value = "-8.78%"
$ws.Cells.Item(43,5).ClearFormats()

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "0.01321"
$ws.Cells.Item(44,4).ClearFormats()
$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,5).Value = "-3.95%"
$ws.Cells.Item(44,5).ClearFormats()

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "0.00006085"
$ws.Cells.Item(45,4).ClearFormats()
$ws.Cells.Item(45,5).NumberFormat = "@"
$ws.Cells.Item(45,5).Value = "1.80%"
$ws.Cells.Item(45,5).ClearFormats()

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "1.917"
$ws.Cells.Item(46,4).ClearFormats()
$ws.Cells.Item(46,5).NumberFormat = "@"
$ws.Cells.Item(46,5).Value = "1.29%"
$ws.Cells.Item(46,5).ClearFormats()

$ws.Cells.Item(47,5).NumberFormat = "@"
$ws.Cells.Item(47,5).Value = "-5.67%"
$ws.Cells.Item(47,5).ClearFormats()
